# Update odds values on the active worksheet (rows 7 and 9) to match
# the latest FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 updates
$ws.Range("G7").Value = 2.52
$ws.Range("I7").Value = 2.75
$ws.Range("S7").Value = 1.8
$ws.Range("T7").Value = 7.1
$ws.Range("U7").Value = 11.75
$ws.Range("V7").Value = 9.75
$ws.Range("W7").Value = 27
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 37
$ws.Range("Z7").Value = 7.9
$ws.Range("AD7").Value = 700
$ws.Range("AE7").Value = 7.9
$ws.Range("AF7").Value = 13.5
$ws.Range("AG7").Value = 10.25
$ws.Range("AH7").Value = 32
$ws.Range("AI7").Value = 25
$ws.Range("AJ7").Value = 37

# Row 9 updates
$ws.Range("G9").Value = 5.5
$ws.Range("H9").Value = 3.85
$ws.Range("I9").Value = 1.55
$ws.Range("S9").Value = 1.98
$ws.Range("V9").Value = 17
$ws.Range("W9").Value = 110
$ws.Range("X9").Value = 50
$ws.Range("Y9").Value = 45
$ws.Range("AA9").Value = 7.8
$ws.Range("AB9").Value = 14.5
$ws.Range("AD9").Value = 400
$ws.Range("AF9").Value = 7.9
$ws.Range("AH9").Value = 11.75
$ws.Range("AJ9").Value = 22
